$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin names and links (columns B and C)
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# Update price and volume columns (D and E) as exact text to preserve formatting
# (avoid Excel auto-converting numeric-looking / percent strings and losing precision/trailing zeros)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.89%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.28%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.034"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.94%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08005"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.95%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.846"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.19%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.111"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.18%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.766"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.74%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9212"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.33%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1270"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.46%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1891"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.05%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09022"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.46%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03446"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.92%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09846"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.39%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001398"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.11%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006278"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.47%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.863"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "7.14%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "14.24%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3412"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.30%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1312"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.81%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.781"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-7.53%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2334"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-10.98%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04344"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.30%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001229"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.79%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004873"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.29%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001298"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.34%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "41.91%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01939"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.70%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05137"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.60%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007502"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.57%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01008"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-10.38%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1352"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.37%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002106"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.12%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009869"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.98%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006162"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.96%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.39%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.85"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.17%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001247"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "4.60%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.39%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.39%"
